$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$r = $ws.Range("C12")
Write-Host ("FormulaHidden: " + $r.FormulaHidden)
$r.FormulaHidden = $true
Write-Host ("After set hidden true: " + $r.FormulaHidden + " Formula=" + $r.Formula + " Value2=" + $r.Value2)
